# Updates odds values on Sheet1 (rows 3, 4, 6, 8) to reflect the latest
# FlashScore odds snapshot, per the "Atualizando o arquivo XLSX" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98

# Row 4
$ws.Range("G4").Value = 1.48
$ws.Range("I4").Value = 7.5
$ws.Range("J4").Value = 2.1
$ws.Range("S4").Value = 2.3
$ws.Range("T4").Value = 1.62
$ws.Range("AD4").Value = 6
$ws.Range("AF4").Value = 9.5
$ws.Range("AP4").Value = 23
$ws.Range("AS4").Value = 67

# Row 6
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.35
$ws.Range("J6").Value = 3.5
$ws.Range("L6").Value = 3
$ws.Range("AC6").Value = 11
$ws.Range("AG6").Value = 23
$ws.Range("AN6").Value = 9
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 9.5

# Row 8
$ws.Range("G8").Value = 3.65
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 1.91
$ws.Range("J8").Value = 4
$ws.Range("K8").Value = 2.18
$ws.Range("L8").Value = 2.5
$ws.Range("N8").Value = 7.9
$ws.Range("T8").Value = 1.98
$ws.Range("W8").Value = 2.77
$ws.Range("X8").Value = 1.39
$ws.Range("Y8").Value = 1.37
$ws.Range("Z8").Value = 2.87
$ws.Range("AB8").Value = 2.1
$ws.Range("AC8").Value = 12.5
$ws.Range("AD8").Value = 21
$ws.Range("AE8").Value = 12
$ws.Range("AF8").Value = 50
$ws.Range("AG8").Value = 30
$ws.Range("AH8").Value = 32
$ws.Range("AI8").Value = 7.9
$ws.Range("AJ8").Value = 6.9
$ws.Range("AL8").Value = 55
$ws.Range("AN8").Value = 8
$ws.Range("AO8").Value = 9.5
$ws.Range("AP8").Value = 8.25
$ws.Range("AQ8").Value = 16.5
$ws.Range("AR8").Value = 14.5
$ws.Range("AS8").Value = 23
